# Auto-generated Excel COM-interop script to apply cell value updates
# as described by the authoritative diff against Phoenix_Profits.xlsx sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# ALC: update 149 cell value(s)
$ws.Range("H28").Value = 1237.8235
$ws.Range("I28").Value = 1278.8
$ws.Range("J28").Value = 930.5
$ws.Range("K28").Value = 1278.8
$ws.Range("L28").Value = 930.5
$ws.Range("M28").Value = -793.8
$ws.Range("N28").Value = -1900.5
$ws.Range("H29").Value = 4119.8
$ws.Range("I29").Value = 849
$ws.Range("J29").Value = 4937.5
$ws.Range("K29").Value = 2547
$ws.Range("L29").Value = 14812.5
$ws.Range("M29").Value = -2266
$ws.Range("N29").Value = -15374.5
$ws.Range("H38").Value = 54843.465
$ws.Range("I38").Value = 73785
$ws.Range("J38").Value = 2754.25
$ws.Range("K38").Value = 221355
$ws.Range("L38").Value = 8262.75
$ws.Range("M38").Value = -220983
$ws.Range("N38").Value = -9006.75
$ws.Range("H40").Value = 5183.2085
$ws.Range("I40").Value = 3362.3333
$ws.Range("J40").Value = 7004.0835
$ws.Range("K40").Value = 3362.3333
$ws.Range("L40").Value = 7004.0835
$ws.Range("M40").Value = -3187.3333
$ws.Range("N40").Value = -7354.0835
$ws.Range("H41").Value = 2318.6428
$ws.Range("I41").Value = 1785.8
$ws.Range("K41").Value = 1785.8
$ws.Range("M41").Value = -1345.8
$ws.Range("H43").Value = 3796.875
$ws.Range("I43").Value = 3695.4
$ws.Range("J43").Value = 3966
$ws.Range("K43").Value = 3695.4
$ws.Range("L43").Value = 3966
$ws.Range("M43").Value = -3626.4
$ws.Range("N43").Value = -4104
$ws.Range("H51").Value = 5003.6
$ws.Range("I51").Value = 5003.6
$ws.Range("K51").Value = 5003.6
$ws.Range("M51").Value = -4519.6
$ws.Range("H58").Value = 146.4
$ws.Range("I58").Value = 146.4
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 439.2
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -289.2
$ws.Range("H64").Value = 8700
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 8700
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 8700
$ws.Range("N64").Value = -9196
$ws.Range("H67").Value = 8700
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 8700
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 8700
$ws.Range("N67").Value = -10416
$ws.Range("H70").Value = 3246.3333
$ws.Range("J70").Value = 3750.75
$ws.Range("L70").Value = 11252.25
$ws.Range("N70").Value = -11792.25
$ws.Range("H73").Value = 3246.3333
$ws.Range("J73").Value = 3750.75
$ws.Range("L73").Value = 11252.25
$ws.Range("N73").Value = -13124.25
$ws.Range("H74").Value = 5042.857
$ws.Range("I74").Value = 5042.857
$ws.Range("K74").Value = 5042.857
$ws.Range("M74").Value = -4106.857
$ws.Range("H77").Value = 5042.857
$ws.Range("I77").Value = 5042.857
$ws.Range("K77").Value = 25214.285
$ws.Range("M77").Value = -20534.285
$ws.Range("H80").Value = 439.1111
$ws.Range("I80").Value = 130
$ws.Range("J80").Value = 593.6667
$ws.Range("K80").Value = 390
$ws.Range("L80").Value = 1781.0001
$ws.Range("M80").Value = 608
$ws.Range("N80").Value = -3777.0001
$ws.Range("H83").Value = 439.1111
$ws.Range("I83").Value = 130
$ws.Range("J83").Value = 593.6667
$ws.Range("K83").Value = 1170
$ws.Range("L83").Value = 5343.0003
$ws.Range("M83").Value = 3822
$ws.Range("N83").Value = -15327.0003
$ws.Range("H86").Value = 2393.889
$ws.Range("J86").Value = 2673.5
$ws.Range("L86").Value = 2673.5
$ws.Range("N86").Value = -4919.5
$ws.Range("H89").Value = 2393.889
$ws.Range("J89").Value = 2673.5
$ws.Range("L89").Value = 13367.5
$ws.Range("N89").Value = -24599.5
$ws.Range("H96").Value = 1138.1305
$ws.Range("I96").Value = 1181.5385
$ws.Range("J96").Value = 1081.7
$ws.Range("K96").Value = 3544.6155
$ws.Range("L96").Value = 3245.1
$ws.Range("M96").Value = -2171.6155
$ws.Range("N96").Value = -5991.1
$ws.Range("H106").Value = 3922.5
$ws.Range("I106").Value = 4341.5
$ws.Range("K106").Value = 4341.5
$ws.Range("M106").Value = -3710.5
$ws.Range("H107").Value = 1053.4736
$ws.Range("I107").Value = 1163.6875
$ws.Range("J107").Value = 465.66666
$ws.Range("K107").Value = 1163.6875
$ws.Range("L107").Value = 465.66666
$ws.Range("M107").Value = 756.3125
$ws.Range("N107").Value = -4305.66666
$ws.Range("H111").Value = 1227
$ws.Range("I111").Value = 1250.3
$ws.Range("J111").Value = 1168.75
$ws.Range("K111").Value = 3750.9
$ws.Range("L111").Value = 3506.25
$ws.Range("M111").Value = -683.8999999999996
$ws.Range("N111").Value = -9640.25
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("H116").Value = 8062.15
$ws.Range("J116").Value = 8911.352999999999
$ws.Range("L116").Value = 8911.352999999999
$ws.Range("N116").Value = -15795.353
$ws.Range("H132").Value = 3117.8157
$ws.Range("I132").Value = 3199.3428
$ws.Range("J132").Value = 2166.6667
$ws.Range("K132").Value = 9598.028399999999
$ws.Range("L132").Value = 6500.000100000001
$ws.Range("M132").Value = -7068.028399999999
$ws.Range("N132").Value = -11560.0001
$ws.Range("H137").Value = 1275.3
$ws.Range("I137").Value = 997.25
$ws.Range("J137").Value = 1460.6666
$ws.Range("K137").Value = 2991.75
$ws.Range("L137").Value = 4381.9998
$ws.Range("M137").Value = -441.75
$ws.Range("N137").Value = -9481.9998
$ws.Range("H138").Value = 2943.3635
$ws.Range("J138").Value = 3068.2856
$ws.Range("L138").Value = 9204.856800000001
$ws.Range("N138").Value = -19484.8568

# ALC: clear 4 cell(s) removed from the sheet
$ws.Range("N58").ClearContents()
$ws.Range("M64").ClearContents()
$ws.Range("M67").ClearContents()
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("ARM")

# ARM: update 36 cell value(s)
$ws.Range("H38").Value = 11019
$ws.Range("I38").Value = 11019
$ws.Range("K38").Value = 11019
$ws.Range("M38").Value = -10552
$ws.Range("H74").Value = 17947.764
$ws.Range("I74").Value = 958.1163
$ws.Range("K74").Value = 958.1163
$ws.Range("M74").Value = -84.11630000000002
$ws.Range("H77").Value = 17947.764
$ws.Range("I77").Value = 958.1163
$ws.Range("K77").Value = 4790.5815
$ws.Range("M77").Value = -422.5815000000002
$ws.Range("H80").Value = 55721.125
$ws.Range("J80").Value = 59395.57
$ws.Range("L80").Value = 59395.57
$ws.Range("N80").Value = -61391.57
$ws.Range("H83").Value = 55721.125
$ws.Range("J83").Value = 59395.57
$ws.Range("L83").Value = 178186.71
$ws.Range("N83").Value = -188170.71
$ws.Range("H110").Value = 1229.909
$ws.Range("I110").Value = 1086.3334
$ws.Range("K110").Value = 1086.3334
$ws.Range("M110").Value = 958.6666
$ws.Range("H115").Value = 66000
$ws.Range("J115").Value = 66000
$ws.Range("L115").Value = 66000
$ws.Range("N115").Value = -69134
$ws.Range("H122").Value = 82966.09
$ws.Range("J122").Value = 91162.7
$ws.Range("L122").Value = 273488.1
$ws.Range("N122").Value = -278388.1
$ws.Range("H132").Value = 2444.6155
$ws.Range("I132").Value = 2394
$ws.Range("K132").Value = 7182
$ws.Range("M132").Value = -4652

$ws = $wb.Worksheets.Item("BSM")

# BSM: update 30 cell value(s)
$ws.Range("H82").Value = 17713.762
$ws.Range("I82").Value = 4935.154
$ws.Range("J82").Value = 38479
$ws.Range("K82").Value = 4935.154
$ws.Range("L82").Value = 38479
$ws.Range("M82").Value = -4552.154
$ws.Range("N82").Value = -39245
$ws.Range("H85").Value = 17713.762
$ws.Range("I85").Value = 4935.154
$ws.Range("J85").Value = 38479
$ws.Range("K85").Value = 4935.154
$ws.Range("L85").Value = 38479
$ws.Range("M85").Value = -3609.154
$ws.Range("N85").Value = -41131
$ws.Range("H99").Value = 2770.7273
$ws.Range("I99").Value = 2735.7896
$ws.Range("K99").Value = 2735.7896
$ws.Range("M99").Value = -1237.7896
$ws.Range("H107").Value = 4477.8
$ws.Range("I107").Value = 4742.722
$ws.Range("K107").Value = 4742.722
$ws.Range("M107").Value = -2822.722
$ws.Range("H114").Value = 69684
$ws.Range("J114").Value = 69684
$ws.Range("L114").Value = 69684
$ws.Range("N114").Value = -78362
$ws.Range("H115").Value = 69684
$ws.Range("J115").Value = 69684
$ws.Range("L115").Value = 69684
$ws.Range("N115").Value = -72818

$ws = $wb.Worksheets.Item("CRP")

# CRP: update 86 cell value(s)
$ws.Range("H14").Value = 5854.25
$ws.Range("I14").Value = 209.5
$ws.Range("K14").Value = 209.5
$ws.Range("M14").Value = -39.5
$ws.Range("H16").Value = 4886.7
$ws.Range("I16").Value = 4812.6665
$ws.Range("J16").Value = 4997.75
$ws.Range("K16").Value = 4812.6665
$ws.Range("L16").Value = 4997.75
$ws.Range("M16").Value = -4525.6665
$ws.Range("N16").Value = -5571.75
$ws.Range("H25").Value = 12326.786
$ws.Range("I25").Value = 11905
$ws.Range("J25").Value = 13381.25
$ws.Range("K25").Value = 11905
$ws.Range("L25").Value = 13381.25
$ws.Range("M25").Value = -11731
$ws.Range("N25").Value = -13729.25
$ws.Range("H86").Value = 10591.211
$ws.Range("I86").Value = 8173.875
$ws.Range("J86").Value = 12349.272
$ws.Range("K86").Value = 8173.875
$ws.Range("L86").Value = 12349.272
$ws.Range("M86").Value = -7050.875
$ws.Range("N86").Value = -14595.272
$ws.Range("H89").Value = 10591.211
$ws.Range("I89").Value = 8173.875
$ws.Range("J89").Value = 12349.272
$ws.Range("K89").Value = 40869.375
$ws.Range("L89").Value = 61746.36
$ws.Range("M89").Value = -35253.375
$ws.Range("N89").Value = -72978.36
$ws.Range("H92").Value = 25632.666
$ws.Range("J92").Value = 25632.666
$ws.Range("L92").Value = 25632.666
$ws.Range("N92").Value = -30624.666
$ws.Range("H96").Value = 13959.444
$ws.Range("J96").Value = 13959.444
$ws.Range("L96").Value = 13959.444
$ws.Range("N96").Value = -19451.444
$ws.Range("H99").Value = 3296.48
$ws.Range("I99").Value = 3428.1333
$ws.Range("J99").Value = 3099
$ws.Range("K99").Value = 3428.1333
$ws.Range("L99").Value = 3099
$ws.Range("M99").Value = -1930.1333
$ws.Range("N99").Value = -6095
$ws.Range("H105").Value = 2045.5
$ws.Range("I105").Value = 2668.25
$ws.Range("J105").Value = 800
$ws.Range("K105").Value = 2668.25
$ws.Range("L105").Value = 800
$ws.Range("M105").Value = -921.25
$ws.Range("N105").Value = -4294
$ws.Range("H107").Value = 52680292
$ws.Range("I107").Value = 76993420
$ws.Range("K107").Value = 76993420
$ws.Range("M107").Value = -76991500
$ws.Range("H113").Value = 4886.7
$ws.Range("I113").Value = 4812.6665
$ws.Range("J113").Value = 4997.75
$ws.Range("K113").Value = 4812.6665
$ws.Range("L113").Value = 4997.75
$ws.Range("M113").Value = -2642.6665
$ws.Range("N113").Value = -9337.75
$ws.Range("H122").Value = 2264.2083
$ws.Range("I122").Value = 1962.7059
$ws.Range("J122").Value = 2996.4285
$ws.Range("K122").Value = 5888.1177
$ws.Range("L122").Value = 8989.2855
$ws.Range("M122").Value = -3438.1177
$ws.Range("N122").Value = -13889.2855
$ws.Range("H126").Value = 3296.48
$ws.Range("I126").Value = 3428.1333
$ws.Range("J126").Value = 3099
$ws.Range("K126").Value = 10284.3999
$ws.Range("L126").Value = 9297
$ws.Range("M126").Value = -7814.3999
$ws.Range("N126").Value = -14237
$ws.Range("H132").Value = 3904.318
$ws.Range("I132").Value = 3103.8667
$ws.Range("J132").Value = 5619.5713
$ws.Range("K132").Value = 9311.6001
$ws.Range("L132").Value = 16858.7139
$ws.Range("M132").Value = -6781.6001
$ws.Range("N132").Value = -21918.7139

$ws = $wb.Worksheets.Item("CUL")

# CUL: update 34 cell value(s)
$ws.Range("H23").Value = 232.71428
$ws.Range("I23").Value = 159.5
$ws.Range("K23").Value = 478.5
$ws.Range("M23").Value = -243.5
$ws.Range("H34").Value = 1096.3889
$ws.Range("J34").Value = 2348.2856
$ws.Range("L34").Value = 7044.8568
$ws.Range("N34").Value = -7212.8568
$ws.Range("H122").Value = 17768.334
$ws.Range("I122").Value = 1305
$ws.Range("J122").Value = 26000
$ws.Range("K122").Value = 11745
$ws.Range("L122").Value = 234000
$ws.Range("M122").Value = -9295
$ws.Range("N122").Value = -238900
$ws.Range("H127").Value = 74998.836
$ws.Range("J127").Value = 74998.836
$ws.Range("L127").Value = 224996.508
$ws.Range("N127").Value = -234916.508
$ws.Range("H132").Value = 1975.3422
$ws.Range("I132").Value = 1352.5834
$ws.Range("J132").Value = 3042.9285
$ws.Range("K132").Value = 12173.2506
$ws.Range("L132").Value = 27386.3565
$ws.Range("M132").Value = -9643.250599999999
$ws.Range("N132").Value = -32446.3565
$ws.Range("H139").Value = 5384.5
$ws.Range("J139").Value = 7373.25
$ws.Range("L139").Value = 22119.75
$ws.Range("N139").Value = -32399.75
$ws.Range("H140").Value = 2595.3333
$ws.Range("J140").Value = 4151.846
$ws.Range("L140").Value = 12455.538
$ws.Range("N140").Value = -22815.538

$ws = $wb.Worksheets.Item("GSM")

# GSM: update 32 cell value(s)
$ws.Range("H21").Value = 20000000
$ws.Range("I21").Value = 20000000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 20000000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -19999827
$ws.Range("H30").Value = 20000000
$ws.Range("I30").Value = 20000000
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 20000000
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -19999895
$ws.Range("H46").Value = 4947.2104
$ws.Range("J46").Value = 39997
$ws.Range("L46").Value = 39997
$ws.Range("N46").Value = -40309
$ws.Range("H64").Value = 59271
$ws.Range("J64").Value = 59271
$ws.Range("L64").Value = 59271
$ws.Range("N64").Value = -59767
$ws.Range("H67").Value = 59271
$ws.Range("J67").Value = 59271
$ws.Range("L67").Value = 59271
$ws.Range("N67").Value = -60987
$ws.Range("H126").Value = 38671.348
$ws.Range("J126").Value = 4717.75
$ws.Range("L126").Value = 14153.25
$ws.Range("N126").Value = -19093.25
$ws.Range("H134").Value = 45581.375
$ws.Range("J134").Value = 45581.375
$ws.Range("L134").Value = 136744.125
$ws.Range("N134").Value = -141814.125

# GSM: clear 2 cell(s) removed from the sheet
$ws.Range("N21").ClearContents()
$ws.Range("N30").ClearContents()

$ws = $wb.Worksheets.Item("LTW")

# LTW: update 48 cell value(s)
$ws.Range("H6").Value = 39785.715
$ws.Range("J6").Value = 39785.715
$ws.Range("L6").Value = 39785.715
$ws.Range("N6").Value = -40009.715
$ws.Range("I16").Value = 882.2
$ws.Range("J16").Value = 9992
$ws.Range("K16").Value = 882.2
$ws.Range("L16").Value = 9992
$ws.Range("M16").Value = -712.2
$ws.Range("N16").Value = -10332
$ws.Range("H22").Value = 4003.2
$ws.Range("J22").Value = 3004.25
$ws.Range("L22").Value = 3004.25
$ws.Range("N22").Value = -3594.25
$ws.Range("H27").Value = 4003.2
$ws.Range("J27").Value = 3004.25
$ws.Range("L27").Value = 3004.25
$ws.Range("N27").Value = -3218.25
$ws.Range("H42").Value = 28028
$ws.Range("I42").Value = 27028
$ws.Range("K42").Value = 27028
$ws.Range("M42").Value = -26465
$ws.Range("H49").Value = 28028
$ws.Range("I49").Value = 27028
$ws.Range("K49").Value = 27028
$ws.Range("M49").Value = -26881
$ws.Range("H82").Value = 1567.4814
$ws.Range("I82").Value = 1550.381
$ws.Range("J82").Value = 1627.3334
$ws.Range("K82").Value = 1550.381
$ws.Range("L82").Value = 1627.3334
$ws.Range("M82").Value = -1189.381
$ws.Range("N82").Value = -2349.3334
$ws.Range("H85").Value = 1567.4814
$ws.Range("I85").Value = 1550.381
$ws.Range("J85").Value = 1627.3334
$ws.Range("K85").Value = 1550.381
$ws.Range("L85").Value = 1627.3334
$ws.Range("M85").Value = -302.3810000000001
$ws.Range("N85").Value = -4123.3334
$ws.Range("H110").Value = 43500
$ws.Range("J110").Value = 43500
$ws.Range("L110").Value = 43500
$ws.Range("N110").Value = -51680
$ws.Range("H136").Value = 50056.473
$ws.Range("I136").Value = 2928.0588
$ws.Range("K136").Value = 8784.1764
$ws.Range("M136").Value = -6234.1764

$ws = $wb.Worksheets.Item("WVR")

# WVR: update 39 cell value(s)
$ws.Range("H18").Value = 10952.1
$ws.Range("J18").Value = 13003.5
$ws.Range("L18").Value = 13003.5
$ws.Range("N18").Value = -13349.5
$ws.Range("H41").Value = 18433.143
$ws.Range("J41").Value = 19006.166
$ws.Range("L41").Value = 19006.166
$ws.Range("N41").Value = -19786.166
$ws.Range("H54").Value = 8500
$ws.Range("J54").Value = 6000
$ws.Range("L54").Value = 6000
$ws.Range("N54").Value = -7040
$ws.Range("H81").Value = 2308.3333
$ws.Range("J81").Value = 10000
$ws.Range("L81").Value = 20000
$ws.Range("N81").Value = -22122
$ws.Range("H84").Value = 2308.3333
$ws.Range("J84").Value = 10000
$ws.Range("L84").Value = 100000
$ws.Range("N84").Value = -110608
$ws.Range("H105").Value = 34999.5
$ws.Range("J105").Value = 34999.5
$ws.Range("L105").Value = 34999.5
$ws.Range("N105").Value = -41987.5
$ws.Range("H111").Value = 66870
$ws.Range("J111").Value = 66870
$ws.Range("L111").Value = 66870
$ws.Range("N111").Value = -75050
$ws.Range("H122").Value = 3008
$ws.Range("I122").Value = 2261.75
$ws.Range("J122").Value = 4998
$ws.Range("K122").Value = 6785.25
$ws.Range("L122").Value = 14994
$ws.Range("M122").Value = -4335.25
$ws.Range("N122").Value = -19894
$ws.Range("H136").Value = 134617230
$ws.Range("I136").Value = 10991123
$ws.Range("K136").Value = 32973369
$ws.Range("M136").Value = -32970819

